$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in cell A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 14.37 = 59399.48 pesos`n✅ 59399.48 pesos = 14.32 = 978.77 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 69.59
$ws2.Range("O10").Value = 4133.61
$ws2.Range("N12").Value = 4148
$ws2.Range("O12").Value = 68.34999999999999
